$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "67.413.52"
Set-TextValue $ws.Range("E2") "  -2.63%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.780.49"
Set-TextValue $ws.Range("E3") "  -0.71%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.07%  "

# Row 5
Set-TextValue $ws.Range("D5") "593.20"
Set-TextValue $ws.Range("E5") "  -1.11%  "

# Row 6
Set-TextValue $ws.Range("D6") "166.12"
Set-TextValue $ws.Range("E6") "  -2.51%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.779.21"
Set-TextValue $ws.Range("E7") "  -0.67%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.08%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.517"
Set-TextValue $ws.Range("E9") "  -1.64%  "

# Row 10
Set-TextValue $ws.Range("E10") "  -2.67%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.35"
Set-TextValue $ws.Range("E11") "  -2.46%  "

# Row 12
Set-TextValue $ws.Range("E12") "  -1.26%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000254"
Set-TextValue $ws.Range("E13") "  -3.70%  "

# Row 14
Set-TextValue $ws.Range("D14") "35.90"
Set-TextValue $ws.Range("E14") "  -2.58%  "

# Row 15
Set-TextValue $ws.Range("D15") "4.415.56"
Set-TextValue $ws.Range("E15") "  -0.70%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.760.06"
Set-TextValue $ws.Range("E16") "  -1.58%  "

# Row 17
Set-TextValue $ws.Range("D17") "67.369.46"
Set-TextValue $ws.Range("E17") "  -2.70%  "

# Row 18
Set-TextValue $ws.Range("E18") "  -0.36%  "

# Row 19
Set-TextValue $ws.Range("E19") "  -0.01%  "

# Row 20
Set-TextValue $ws.Range("D20") "6.97"
Set-TextValue $ws.Range("E20") "  -1.84%  "

# Row 21
Set-TextValue $ws.Range("D21") "10.18"
Set-TextValue $ws.Range("E21") "  -7.40%  "

# Row 22
Set-TextValue $ws.Range("D22") "457.74"
Set-TextValue $ws.Range("E22") "  -2.99%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.698"
Set-TextValue $ws.Range("E23") "  -1.53%  "

# Row 24
Set-TextValue $ws.Range("E24") "  +1.55%  "

# Row 25
Set-TextValue $ws.Range("D25") "83.41"
Set-TextValue $ws.Range("E25") "  -1.86%  "

# Row 26
Set-TextValue $ws.Range("B26") "Fetch.AI"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D26") "2.13"
Set-TextValue $ws.Range("E26") "  -5.12%  "

# Row 27
Set-TextValue $ws.Range("B27") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D27") "11.82"
Set-TextValue $ws.Range("E27") "  -3.21%  "

# Row 28
Set-TextValue $ws.Range("E28") "  +0.06%  "

# Row 29
Set-TextValue $ws.Range("D29") "9.97"
Set-TextValue $ws.Range("E29") "  -3.03%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -1.77%  "

# Row 31
Set-TextValue $ws.Range("D31") "29.78"
Set-TextValue $ws.Range("E31") "  -1.82%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -2.61%  "

# Row 33
Set-TextValue $ws.Range("E33") "  -4.05%  "

# Row 34
Set-TextValue $ws.Range("E34") "  -3.11%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  +0.02%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.734.22"
Set-TextValue $ws.Range("E36") "  -0.77%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -2.62%  "

# Row 38
Set-TextValue $ws.Range("E38") "  -6.83%  "

# Row 39
Set-TextValue $ws.Range("E39") "  -1.47%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.994"
Set-TextValue $ws.Range("E40") "  -1.88%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -3.05%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.999"
Set-TextValue $ws.Range("E42") "  -0.13%  "

# Row 43
Set-TextValue $ws.Range("E43") "  +0.00%  "

# Row 44
Set-TextValue $ws.Range("D44") "44.04"
Set-TextValue $ws.Range("E44") "  -0.21%  "

# Row 45
Set-TextValue $ws.Range("D45") "47.02"
Set-TextValue $ws.Range("E45") "  +1.85%  "

# Row 46
Set-TextValue $ws.Range("E46") "  -4.45%  "

# Row 47
Set-TextValue $ws.Range("E47") "  -3.71%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +1.20%  "

# Row 49
Set-TextValue $ws.Range("D49") "392.80"
Set-TextValue $ws.Range("E49") "  -2.19%  "

# Row 50
Set-TextValue $ws.Range("E50") "  -8.42%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.754.58"
Set-TextValue $ws.Range("E51") "  +1.83%  "
